$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the existing row 126, pushing it down to row 127
$ws.Rows.Item(126).Insert()

# Populate the new row 126 with the updated weekly data
$ws.Cells.Item(126, 1).Value = 8
$ws.Cells.Item(126, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(126, 3).Value = "Coquimbo"
$ws.Cells.Item(126, 4).Value = 44628
$ws.Cells.Item(126, 5).Value = 4
$ws.Cells.Item(126, 6).Value = 100112044
$ws.Cells.Item(126, 7).Value = "Perejil"
$ws.Cells.Item(126, 8).Value = "Sin especificar"
$ws.Cells.Item(126, 9).Value = "Primera"
$ws.Cells.Item(126, 10).Value = 2400
$ws.Cells.Item(126, 11).Value = 2500
$ws.Cells.Item(126, 12).Value = 3000
$ws.Cells.Item(126, 13).Value = 2750
$ws.Cells.Item(126, 14).Value = '$/atado 1 a 1,5 kilos'
$ws.Cells.Item(126, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(126, 16).Value = 1833
$ws.Cells.Item(126, 17).Value = 1.5
$ws.Cells.Item(126, 18).Value = "Hortaliza"

# Ensure the date cell keeps the date/time number format used by column D
$ws.Cells.Item(126, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
